$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$win.ScrollRow = 152
$win.ScrollColumn = 1
$sr = $win.ScrollRow
Write-Host "before save scrollrow=$sr"
